# Add a "Type" / "BASIC_TEXT_TYPE" row to the basic-text block (A79:B82 area)
# so the badge table gains a background/type column, allowing customisation
# of different badge sizes.
#
# This mirrors inserting a new row 80 in Excel: everything from the old
# row 80 ("Location") downward shifts down by one row, formulas/ranges and
# the sheet dimension auto-adjust, and the new row is populated with the
# "Type" label pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 80 (shifts Location/Description/... + scenario rows down by one)
$ws.Rows.Item(80).Insert()

# Populate the new row with the Type label pair (same shape as the other
# BASIC_TEXT_* rows: col A holds the key, col B the display text, and
# C:E hold the "XXXX" placeholder used throughout this table)
$ws.Range("A80").Value = "BASIC_TEXT_TYPE"
$ws.Range("B80").Value = "Type"
$ws.Range("C80").Value = "XXXX"
$ws.Range("D80").Value = "XXXX"
$ws.Range("E80").Value = "XXXX"

# Move the selection/scroll position to roughly where the edit happened
$aw = $excel.ActiveWindow
$aw.ScrollRow = 65
$aw.ScrollColumn = 3
$ws.Range("E80").Select()
